$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data refresh.
# D/E columns hold plain text in the source data (e.g. "64.296.21" is not
# a valid number, and values like "569.85" must stay text, not become a
# float) so NumberFormat is forced to Text ("@") before writing each value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.296.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.401.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.45%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.96"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +9.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.405.35"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.36%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.989.03"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.334.90"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.364.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.80"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.99"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.45%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.72"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.96%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.33"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.16%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.12%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.05"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.63%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.34%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.41"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.91%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0758"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.32%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.862.97"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.38"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.19%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.73"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.79"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.766"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.111"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.60%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.36%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.51%  "

